$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows 8-15 (campaign list shrank to 6 data rows)
$ws.Rows("8:15").Delete()

# Overwrite the remaining data rows (A2:A7) with the new campaign names.
# The header in A1 ("Campaign Name") is left untouched.
$ws.Cells.Item(2, 1).Value = "1LM/HAR/701_Chevy North East LMA_HARTFORD-NEW HAVEN_LMA_BCN_NER_GNM_2018"
$ws.Cells.Item(3, 1).Value = "UMM/UM/683_Maven_Maven_Peer_BCN_USA_GNM_2018"
$ws.Cells.Item(4, 1).Value = "BUN/BD/686_Buick Division_Buick Brand_3Q'18 OLV Digital Video_BLA_USA_GNM_2018_Buick_SPORTS_STREAM"
$ws.Cells.Item(5, 1).Value = "BUN/BD/686_Buick Division_Buick Nat'l Full Line_2018 CBS College Football Live Streaming_BPR_USA_GNM_2018"
$ws.Cells.Item(6, 1).Value = "CVN/CEA/683_Chevrolet Division_Chevy Equinox AA_Equinox AA_DTU Amplification_AWR_USA_AAM_2018"
$ws.Cells.Item(7, 1).Value = "CAX/DM/683_Cadillac Division_Cadillac Brand_Variety_BLA_USA_GNM_2018							"
